# Weekly data refresh: insert the newest price-report row for
# "Terminal La Palmera de La Serena - Papa" at the top of the dated
# block (row 595), pushing the existing rows (595-628) down by one
# (596-629). This mirrors how Excel's Insert Row shifts the sheet and
# grows the used range from A1:R628 to A1:R629.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 595; everything from 595 downward
# shifts down one row (595->596, ..., 628->629).
$ws.Rows.Item(595).Insert()

# Populate the newly inserted row 595 with the new weekly report data.
$ws.Cells.Item(595, 1).Value  = 8
$ws.Cells.Item(595, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(595, 3).Value  = "Coquimbo"
$ws.Cells.Item(595, 4).Value  = 45041
$ws.Cells.Item(595, 5).Value  = 4
$ws.Cells.Item(595, 6).Value  = 100114001
$ws.Cells.Item(595, 7).Value  = "Papa"
$ws.Cells.Item(595, 8).Value  = "Asterix"
$ws.Cells.Item(595, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(595, 10).Value = 2000
$ws.Cells.Item(595, 11).Value = 11500
$ws.Cells.Item(595, 12).Value = 12000
$ws.Cells.Item(595, 13).Value = 11750
$ws.Cells.Item(595, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(595, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(595, 16).Value = 470
$ws.Cells.Item(595, 17).Value = 25
$ws.Cells.Item(595, 18).Value = "Hortaliza"
